$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 44946
$ws.Range("C8").Value = 44946
$ws.Range("D8").Value = 44946
$ws.Range("E8").Value = 44946

$ws.Range("E8").Select()
